$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 31
$ws.Range("I2").Value = 26804
$ws.Range("J2").Value = 2012
$ws.Range("K2").Value = "5fab3d3029d511f60544d7a58d047e5c"
$ws.Range("L2").Value = "GABE2**081978"
$ws.Range("M2").Value = 41240
$ws.Range("N2").Value = 41053
$ws.Range("O2").Value = 41238
$ws.Range("P2").Value = 146
$ws.Range("Q2").Value = 187
$ws.Range("R2").Value = 185
$ws.Range("S2").Value = "Derivación"
$ws.Range("T2").Value = "Derivación"
$ws.Range("U2").Value = "PG-PAB"
$ws.Range("V2").Value = 15484
$ws.Range("W2").Value = 15669
$ws.Range("X2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 31
$ws.Range("I3").Value = 31405
$ws.Range("J3").Value = 2012
$ws.Range("K3").Value = "5fab3d3029d511f60544d7a58d047e5c"
$ws.Range("L3").Value = "GABE2**081978"
$ws.Range("M3").Value = 41260
$ws.Range("N3").Value = 41207
$ws.Range("O3").Value = 41238
$ws.Range("P3").Value = 142
$ws.Range("Q3").Value = 53
$ws.Range("R3").Value = 31
$ws.Range("S3").Value = "Alta Terapéutica"
$ws.Range("T3").Value = "Alta Terapéutica"
$ws.Range("U3").Value = "PG-PR"
$ws.Range("V3").Value = 15638
$ws.Range("W3").Value = 15669
$ws.Range("X3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 6
$ws.Range("I4").Value = 36873
$ws.Range("J4").Value = 2013
$ws.Range("K4").Value = "953090fd0a141e1cb79c86c7b209009f"
$ws.Range("L4").Value = "JOPA1**031985"
$ws.Range("M4").Value = 41449
$ws.Range("N4").Value = 41295
$ws.Range("O4").Value = 41437
$ws.Range("P4").Value = 155
$ws.Range("Q4").Value = 154
$ws.Range("R4").Value = 142
$ws.Range("S4").Value = "Abandono"
$ws.Range("T4").Value = "Abandono"
$ws.Range("U4").Value = "PG-PAB"
$ws.Range("V4").Value = 15726
$ws.Range("W4").Value = 15868
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 6
$ws.Range("I5").Value = 42419
$ws.Range("J5").Value = 2013
$ws.Range("K5").Value = "953090fd0a141e1cb79c86c7b209009f"
$ws.Range("L5").Value = "JOPA1**031985"
$ws.Range("M5").Value = 41607
$ws.Range("N5").Value = 41431
$ws.Range("O5").Value = 41607
$ws.Range("P5").Value = 254
$ws.Range("Q5").Value = 176
$ws.Range("R5").Value = 176
$ws.Range("S5").Value = "Alta Admnistrativa"
$ws.Range("T5").Value = "Alta Admnistrativa"
$ws.Range("U5").Value = "PG-PAI"
$ws.Range("V5").Value = 15862
$ws.Range("W5").Value = 16038
$ws.Range("C6").Value = 1
$ws.Range("G6").Value = 71
$ws.Range("I6").Value = 157931
$ws.Range("J6").Value = 2019
$ws.Range("K6").Value = "a92ab8f10fa05d0d9dcb5855c0ec0092"
$ws.Range("L6").Value = "RITR1**111984"
$ws.Range("N6").Value = 43579
$ws.Range("O6").Value = 43669
$ws.Range("P6").Value = 291
$ws.Range("Q6").Value = 203
$ws.Range("R6").Value = 90
$ws.Range("V6").Value = 18010
$ws.Range("W6").Value = 18100
$ws.Range("C7").Value = 1
$ws.Range("G7").Value = 71
$ws.Range("I7").Value = 156785
$ws.Range("J7").Value = 2019
$ws.Range("K7").Value = "a92ab8f10fa05d0d9dcb5855c0ec0092"
$ws.Range("L7").Value = "RITR1**111984"
$ws.Range("M7").Value = 43607
$ws.Range("N7").Value = 43598
$ws.Range("O7").Value = 43607
$ws.Range("P7").Value = 291
$ws.Range("Q7").Value = 9
$ws.Range("R7").Value = 9
$ws.Range("S7").Value = "Abandono"
$ws.Range("T7").Value = "Abandono"
$ws.Range("U7").Value = "M-PAI"
$ws.Range("V7").Value = 18029
$ws.Range("W7").Value = 18038
$ws.Range("C8").Value = 0
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = 90
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 111197
$ws.Range("J8").Value = 2017
$ws.Range("K8").Value = "144b8d70d7ea1b9ea70d2ef7543520b2"
$ws.Range("L8").Value = "OSBE1**061978"
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = 42720
$ws.Range("O8").Value = 42856
$ws.Range("P8").Value = 205
$ws.Range("Q8").Value = 1054
$ws.Range("R8").Value = 136
$ws.Range("S8").Value = $null
$ws.Range("T8").Value = $null
$ws.Range("V8").Value = 17151
$ws.Range("W8").Value = 17287
$ws.Range("X8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("F9").Value = $null
$ws.Range("G9").Value = 90
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 112951
$ws.Range("J9").Value = 2017
$ws.Range("K9").Value = "144b8d70d7ea1b9ea70d2ef7543520b2"
$ws.Range("L9").Value = "OSBE1**061978"
$ws.Range("M9").Value = 43038
$ws.Range("N9").Value = 42766
$ws.Range("O9").Value = 43038
$ws.Range("P9").Value = 682
$ws.Range("Q9").Value = 272
$ws.Range("R9").Value = 272
$ws.Range("S9").Value = "Alta Terapéutica"
$ws.Range("T9").Value = "Alta Terapéutica"
$ws.Range("U9").Value = "PG-PR"
$ws.Range("V9").Value = 17197
$ws.Range("W9").Value = 17469
$ws.Range("X9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("G10").Value = 23
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 109317
$ws.Range("J10").Value = 2017
$ws.Range("K10").Value = "1baff032eb17af74d98d63542c87423a"
$ws.Range("L10").Value = "CEME1**081985"
$ws.Range("N10").Value = 42626
$ws.Range("O10").Value = 42785
$ws.Range("P10").Value = 205
$ws.Range("Q10").Value = 1148
$ws.Range("R10").Value = 159
$ws.Range("U10").Value = "PG-PAI"
$ws.Range("V10").Value = 17057
$ws.Range("W10").Value = 17216
$ws.Range("X10").Value = 1
$ws.Range("C11").Value = 0
$ws.Range("G11").Value = 23
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 112950
$ws.Range("J11").Value = 2017
$ws.Range("K11").Value = "1baff032eb17af74d98d63542c87423a"
$ws.Range("L11").Value = "CEME1**081985"
$ws.Range("M11").Value = 42775
$ws.Range("N11").Value = 42762
$ws.Range("O11").Value = 42775
$ws.Range("P11").Value = 682
$ws.Range("Q11").Value = 13
$ws.Range("R11").Value = 13
$ws.Range("S11").Value = "Alta Admnistrativa"
$ws.Range("T11").Value = "Alta Admnistrativa"
$ws.Range("V11").Value = 17193
$ws.Range("W11").Value = 17206
$ws.Range("X11").Value = 1
$ws.Range("C12").Value = $null
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 511
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 86673
$ws.Range("J12").Value = 2016
$ws.Range("K12").Value = "380d9b27915ca8873d06f71d4d74030e"
$ws.Range("L12").Value = "JUCA1**071955"
$ws.Range("M12").Value = 42458
$ws.Range("N12").Value = 41939
$ws.Range("O12").Value = 42458
$ws.Range("P12").Value = 105
$ws.Range("Q12").Value = 519
$ws.Range("R12").Value = 519
$ws.Range("S12").Value = "Alta Terapéutica"
$ws.Range("T12").Value = "Alta Terapéutica"
$ws.Range("U12").Value = "PG-PAI"
$ws.Range("V12").Value = 16370
$ws.Range("W12").Value = 16889
$ws.Range("X12").Value = 1
$ws.Range("C13").Value = $null
$ws.Range("E13").Value = "Si"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 511
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 64530
$ws.Range("K13").Value = "380d9b27915ca8873d06f71d4d74030e"
$ws.Range("L13").Value = "JUCA1**071955"
$ws.Range("M13").Value = 41948
$ws.Range("N13").Value = 41947
$ws.Range("O13").Value = 41948
$ws.Range("P13").Value = $null
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 1
$ws.Range("U13").Value = "CALLE"
$ws.Range("V13").Value = 16378
$ws.Range("W13").Value = 16379
$ws.Range("X13").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = $null
$ws.Range("G14").Value = 88
$ws.Range("I14").Value = 89750
$ws.Range("J14").Value = 2016
$ws.Range("K14").Value = "a91ebc49e725f0638be44c6e17445adb"
$ws.Range("L14").Value = "BACA1**071992"
$ws.Range("M14").Value = $null
$ws.Range("N14").Value = 42270
$ws.Range("O14").Value = 42366
$ws.Range("P14").Value = 167
$ws.Range("Q14").Value = 1504
$ws.Range("R14").Value = 96
$ws.Range("S14").Value = $null
$ws.Range("T14").Value = $null
$ws.Range("U14").Value = "PG-PR"
$ws.Range("V14").Value = 16701
$ws.Range("W14").Value = 16797
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = $null
$ws.Range("G15").Value = 88
$ws.Range("I15").Value = 94274
$ws.Range("J15").Value = 2016
$ws.Range("K15").Value = "a91ebc49e725f0638be44c6e17445adb"
$ws.Range("L15").Value = "BACA1**071992"
$ws.Range("M15").Value = 42475
$ws.Range("N15").Value = 42278
$ws.Range("O15").Value = 42475
$ws.Range("P15").Value = 167
$ws.Range("Q15").Value = 197
$ws.Range("R15").Value = 197
$ws.Range("S15").Value = "Abandono"
$ws.Range("T15").Value = "Abandono"
$ws.Range("U15").Value = "PG-PR"
$ws.Range("V15").Value = 16709
$ws.Range("W15").Value = 16906
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = $null
$ws.Range("G16").Value = 214
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 65674
$ws.Range("J16").Value = 2015
$ws.Range("K16").Value = "137e8525aa3f79235fa8ad90913fdcbe"
$ws.Range("L16").Value = "FRCA1**071985"
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = 41449
$ws.Range("O16").Value = 41672
$ws.Range("P16").Value = 139
$ws.Range("Q16").Value = 2325
$ws.Range("R16").Value = 223
$ws.Range("S16").Value = $null
$ws.Range("T16").Value = $null
$ws.Range("V16").Value = 15880
$ws.Range("W16").Value = 16103
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = "No"
$ws.Range("F17").Value = $null
$ws.Range("G17").Value = 214
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 49659
$ws.Range("J17").Value = 2014
$ws.Range("K17").Value = "137e8525aa3f79235fa8ad90913fdcbe"
$ws.Range("L17").Value = "FAMU2**071985"
$ws.Range("M17").Value = 41791
$ws.Range("N17").Value = 41458
$ws.Range("O17").Value = 41791
$ws.Range("P17").Value = 139
$ws.Range("Q17").Value = 333
$ws.Range("R17").Value = 333
$ws.Range("S17").Value = "Alta Admnistrativa"
$ws.Range("T17").Value = "Alta Admnistrativa"
$ws.Range("U17").Value = "PG-PAI"
$ws.Range("V17").Value = 15889
$ws.Range("W17").Value = 16222

Write-Host "done"